# Apply the edits described by the diff to the g11.8a workbook.
# Years move from 2022 -> 2023, the region table is refreshed with a new
# set of percentages/ranks and a new list of region labels, and the label
# column (A) and numeric columns (B/C) pick up new number formats /
# borders / fonts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: year 2022 -> 2023 -----------------------------------------------
$ws.Range("B1").Value = 2023

# --- Column A labels (rows 6-10) ------------------------------------------
# Row 5 ("PA") keeps its text; the others get new region codes.  The order
# below matches the order the strings appear in the target shared-string
# table once the now-unused "AP"/"CE" entries are dropped.
$ws.Range("A6").Value = "MA"
$ws.Range("A7").Value = "MT"
$ws.Range("A8").Value = "TO"
$ws.Range("A9").Value = "DF"
$ws.Range("A10").Value = "RR"

# --- Column B values (rows 2-10), now expressed as fractions --------------
$ws.Range("B2").Value = -0.98946394060239751
$ws.Range("B3").Value = -0.92350451060959049
$ws.Range("B4").Value = -1.1528940036622264
$ws.Range("B5").Value = -0.22005518309996197
$ws.Range("B6").Value = -0.18400518517543704
$ws.Range("B7").Value = -0.14318086433921523
$ws.Range("B8").Value = -0.026650710527051397
$ws.Range("B9").Value = 0.064945193163240517
$ws.Range("B10").Value = 0.19726082681291643

# --- Column C values -------------------------------------------------------
$ws.Range("C4").Value = 13
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 3
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 1

# --- Number formats ---------------------------------------------------------
$ws.Range("B2:B10").NumberFormat = "0.0%"
$ws.Range("C2").NumberFormat = "0.0"
$ws.Range("C3").NumberFormat = "0.0"
$ws.Range("C4:C10").NumberFormat = "0"

# --- Font + border + wrap for the label column (applied cell by cell so ---
# --- every row in a block converges on the same cell style) ---------------
for ($r = 2; $r -le 3; $r++) {
    $c = $ws.Cells.Item($r, 1)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.ColorIndex = 8
    $c.Borders.Item(7).ColorIndex = 22
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(10).ColorIndex = 22
    $c.Borders.Item(10).LineStyle = 1
    $c.WrapText = $true
}

for ($r = 4; $r -le 10; $r++) {
    $c = $ws.Cells.Item($r, 1)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.ColorIndex = 8
    $c.Borders.Item(7).ColorIndex = 22
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(8).ColorIndex = 22
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(9).ColorIndex = 22
    $c.Borders.Item(9).LineStyle = 1
    $c.Borders.Item(10).ColorIndex = 22
    $c.Borders.Item(10).LineStyle = 1
    $c.WrapText = $true
}

# --- Selection --------------------------------------------------------------
$ws.Range("A2:C10").Select()

Write-Output "edit complete"
